# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps (and the rolled-up "Latest HO Xliff Generate Date" on the
# Overview sheet) for the 57254883-6acc-4743-b372-808c31626313.md file,
# reflecting a newly generated handback/handoff xliff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is 57254883-6acc-4743-b372-808c31626313.md ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G3").Value = "2016-12-15 04:12:02"

# --- zh-cn sheet: row 3 is 57254883-6acc-4743-b372-808c31626313.md ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("H3").Value = "2016-12-15 04:11:48"
$ws.Range("L3").Value = "2016-12-15 04:12:41"

# --- de-de sheet: row 3 is 57254883-6acc-4743-b372-808c31626313.md ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("H3").Value = "2016-12-15 04:12:02"
$ws.Range("L3").Value = "2016-12-15 04:12:59"
